$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 13: replace the inspector entry. The old "FERNANDEZ MAURICIO LORENZO"
# row (row 13) is removed and the data that used to live in row 14
# ("SDFSDFDS") is moved up into row 13. Row 14 becomes blank.
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "SDFSDFDS"
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(13, 10).Value = 29
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 8
$ws.Cells.Item(13, 15).Value = 3
$ws.Cells.Item(13, 48).Value = 9.0

# Row 14: clear every column (A:AW) so the row becomes empty.
$ws.Range("A14:AW14").ClearContents()
